$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    # Coinranking "Price" values are stored as literal text (not numbers) in
    # the source sheet (e.g. thousands-grouped "26.665.51", or plain decimals
    # like "1.00" that Excel would otherwise auto-convert to a Number).
    # Force text storage, write the value, then drop the temporary "@" text
    # format so the cell keeps its original (default) style.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).ClearFormats()
}

# --- Update price (D) and 1h-volume-change (E) columns for rows 2-47 ---
Set-TextValue $ws "D2" "26.665.51"
$ws.Range("E2").Value = "  +0.05%  "
Set-TextValue $ws "D3" "1.597.88"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws "D5" "211.52"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue $ws "D8" "0.0618"
$ws.Range("E8").Value = "  +0.12%  "
Set-TextValue $ws "D9" "0.248"
$ws.Range("E9").Value = "  +0.76%  "
Set-TextValue $ws "D10" "19.57"
$ws.Range("E10").Value = "  -0.48%  "
Set-TextValue $ws "D11" "0.0841"
$ws.Range("E11").Value = "  +0.36%  "
Set-TextValue $ws "D12" "1.822.49"
$ws.Range("E12").Value = "  +0.11%  "
Set-TextValue $ws "D13" "1.646.58"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  +0.49%  "
Set-TextValue $ws "D16" "65.17"
$ws.Range("E16").Value = "  +0.39%  "
Set-TextValue $ws "D17" "26.655.62"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +1.34%  "
Set-TextValue $ws "D19" "209.39"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +0.17%  "
Set-TextValue $ws "D21" "7.02"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E23").Value = "  +1.74%  "
Set-TextValue $ws "D24" "8.97"
$ws.Range("E24").Value = "  +0.72%  "
Set-TextValue $ws "D25" "144.29"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  -0.24%  "
Set-TextValue $ws "D29" "15.28"
$ws.Range("E29").Value = "  -0.03%  "
Set-TextValue $ws "D30" "0.0516"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +1.58%  "
Set-TextValue $ws "D34" "1.289.06"
$ws.Range("E34").Value = "  -0.73%  "
Set-TextValue $ws "D35" "0.618"
$ws.Range("E35").Value = "  -7.03%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("E38").Value = "  -0.61%  "
Set-TextValue $ws "D39" "0.834"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  +19.38%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  -0.72%  "
Set-TextValue $ws "D44" "63.54"
$ws.Range("E44").Value = "  -0.43%  "
Set-TextValue $ws "D45" "1.735.58"
$ws.Range("E45").Value = "  +0.11%  "
Set-TextValue $ws "D46" "90.76"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("E47").Value = "  -3.05%  "

# --- Rows 48-51: BabyDogeCoin row removed; Algorand/Cronos/USDD rows each
# shift up by one; a new "EnergySwap" row is appended at the bottom ---
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D48" "0.102"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D49" "0.0508"
$ws.Range("E49").Value = "  +0.85%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws "D50" "1.00"
$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D51" "7.40"
$ws.Range("E51").Value = "  -0.94%  "
